# Added Panel Accessories Test Data For Spain/Turkey/Hungary market
#
# For the Spain, Turkey and Hungary worksheets, insert two new accessory
# rows ("MX-BBX" and "MX-DPBX") right after the existing "PR1D2-Unmonitored"
# row (and before the trailing "Wg"/"Accessories" marker rows), matching
# the layout already used on other market sheets (e.g. Portugal).

$wb = $excel.ActiveWorkbook

$targetSheets = @("Spain", "Turkey", "Hungary")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Push the old rows 11 ("Wg") and 12 ("Accessories") down to 13/14,
    # opening up rows 11/12 for the two new accessory entries.
    $ws.Rows("11:12").Insert()

    # Copy the formatting of the row above (style index 3, bordered cell)
    # onto the two freshly inserted rows so they match the rest of the list.
    $ws.Range("A10").Copy()
    $ws.Range("A11:A12").PasteSpecial(-4122)

    $ws.Range("A11").Value = "MX-BBX"
    $ws.Range("A12").Value = "MX-DPBX"

    # Leave the newly added rows selected on this sheet.
    $ws.Range("A11:A12").Select()
}

# The Portugal sheet keeps its existing selection rectangle over its own
# MX-BBX/MX-DPBX rows, but is no longer the active tab.
$wsPortugal = $wb.Worksheets.Item("Portugal")
$wsPortugal.Range("A9:A10").Select()

# Hungary becomes the active sheet/tab after this round of edits.
$wsHungary = $wb.Worksheets.Item("Hungary")
$wsHungary.Activate()
